$d = $word.ActiveDocument

# 1) Merge the three split hyperlink runs ("Basic Syntax | Markdo" + "w" + "n Guide")
#    into a single run with the same visible text "Basic Syntax | Markdown Guide".
#    Find/Replace (wdReplaceAll=2) forces Word to rewrite the run, coalescing the
#    three runs that previously made up this text into one.
$d.Content.Find.Execute("Basic Syntax | Markdown Guide", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Basic Syntax | Markdown Guide", 2) | Out-Null

# 2) Remove the existing _GoBack bookmark (currently sitting at the end of the
#    "- Na pasta q eu quero criar o repositório" paragraph).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 3) Insert a new "git status" paragraph right after "git add ." (before "git commit -m ...").
$gitAddPara = $d.Paragraphs(10)
$gitAddPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(11)
$newPara.Range.Text = "git status"

# 4) Re-create the _GoBack bookmark, now collapsed at the end of the new
#    "git status" paragraph's text (matching its new location in the target doc).
$newParaRange = $d.Paragraphs(11).Range
$bmRange = $d.Range($newParaRange.Start, $newParaRange.Start + 10)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
